# Auto-generated Excel COM-interop script
# Applies the cell-value updates described by the commit diff
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 10000000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 10000000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 10000000
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -10000340
# Row 62
$ws.Range("H62").Value = 55557220
$ws.Range("I62").Value = 66668068
$ws.Range("J62").Value = 2995.3333
$ws.Range("K62").Value = 66668068
$ws.Range("L62").Value = 2995.3333
$ws.Range("M62").Value = -66667444
$ws.Range("N62").Value = -4243.3333
# Row 65
$ws.Range("H65").Value = 55557220
$ws.Range("I65").Value = 66668068
$ws.Range("J65").Value = 2995.3333
$ws.Range("K65").Value = 333340340
$ws.Range("L65").Value = 14976.6665
$ws.Range("M65").Value = -333337220
$ws.Range("N65").Value = -21216.6665
# Row 107
$ws.Range("H107").Value = 580.4828
$ws.Range("I107").Value = 580.4828
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 580.4828
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1339.5172
$ws.Range("N107").ClearContents()
# Row 111
$ws.Range("H111").Value = 1354.2941
$ws.Range("I111").Value = 1606.5555
$ws.Range("J111").Value = 1070.5
$ws.Range("K111").Value = 4819.666499999999
$ws.Range("L111").Value = 3211.5
$ws.Range("M111").Value = -1752.666499999999
$ws.Range("N111").Value = -9345.5
# Row 116
$ws.Range("H116").Value = 209181.86
$ws.Range("I116").Value = 5345.4546
$ws.Range("J116").Value = 268187.12
$ws.Range("K116").Value = 5345.4546
$ws.Range("L116").Value = 268187.12
$ws.Range("M116").Value = -1903.4546
$ws.Range("N116").Value = -275071.12
# Row 127
$ws.Range("H127").Value = 2444.8235
$ws.Range("I127").Value = 973.125
$ws.Range("J127").Value = 3753
$ws.Range("K127").Value = 2919.375
$ws.Range("L127").Value = 11259
$ws.Range("M127").Value = 2040.625
$ws.Range("N127").Value = -21179
# Row 137
$ws.Range("H137").Value = 1209.875
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 1209.875
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 3629.625
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -8729.625
# Row 138
$ws.Range("H138").Value = 1511.33
$ws.Range("I138").Value = 792.8333
$ws.Range("J138").Value = 2174.5576
$ws.Range("K138").Value = 2378.4999
$ws.Range("L138").Value = 6523.6728
$ws.Range("M138").Value = 2761.5001
$ws.Range("N138").Value = -16803.6728
# Row 141
$ws.Range("H141").Value = 2440.4827
$ws.Range("I141").Value = 2055.8286
$ws.Range("J141").Value = 4024.353
$ws.Range("K141").Value = 6167.485799999999
$ws.Range("L141").Value = 12073.059
$ws.Range("M141").Value = -987.4857999999995
$ws.Range("N141").Value = -22433.059

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2260
$ws.Range("I2").Value = 2150
$ws.Range("J2").Value = 2333.3333
$ws.Range("K2").Value = 2150
$ws.Range("L2").Value = 2333.3333
$ws.Range("M2").Value = -2037
$ws.Range("N2").Value = -2559.3333
# Row 32
$ws.Range("H32").Value = 5354.099
$ws.Range("I32").Value = 4579.6533
$ws.Range("J32").Value = 8984.3125
$ws.Range("K32").Value = 4579.6533
$ws.Range("L32").Value = 8984.3125
$ws.Range("M32").Value = -4292.6533
$ws.Range("N32").Value = -9558.3125
# Row 61
$ws.Range("H61").Value = 4245.56
$ws.Range("I61").Value = 4245.56
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4245.56
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4033.56
$ws.Range("N61").ClearContents()
# Row 74
$ws.Range("H74").Value = 785.93335
$ws.Range("I74").Value = 697.44446
$ws.Range("J74").Value = 918.6667
$ws.Range("K74").Value = 697.44446
$ws.Range("L74").Value = 918.6667
$ws.Range("M74").Value = 176.55554
$ws.Range("N74").Value = -2666.6667
# Row 77
$ws.Range("H77").Value = 785.93335
$ws.Range("I77").Value = 697.44446
$ws.Range("J77").Value = 918.6667
$ws.Range("K77").Value = 3487.2223
$ws.Range("L77").Value = 4593.3335
$ws.Range("M77").Value = 880.7776999999996
$ws.Range("N77").Value = -13329.3335
# Row 116
$ws.Range("H116").Value = 2260
$ws.Range("I116").Value = 2150
$ws.Range("J116").Value = 2333.3333
$ws.Range("K116").Value = 2150
$ws.Range("L116").Value = 2333.3333
$ws.Range("M116").Value = 144
$ws.Range("N116").Value = -6921.3333
# Row 132
$ws.Range("H132").Value = 3499.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -15558.5
# Row 136
$ws.Range("H136").Value = 4245.56
$ws.Range("I136").Value = 4245.56
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12736.68
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10186.68
$ws.Range("N136").ClearContents()

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2260
$ws.Range("I3").Value = 2150
$ws.Range("J3").Value = 2333.3333
$ws.Range("K3").Value = 2150
$ws.Range("L3").Value = 2333.3333
$ws.Range("M3").Value = -2036
$ws.Range("N3").Value = -2561.3333
# Row 107
$ws.Range("H107").Value = 47620520
$ws.Range("I107").Value = 52633052
$ws.Range("J107").Value = 1450
$ws.Range("K107").Value = 52633052
$ws.Range("L107").Value = 1450
$ws.Range("M107").Value = -52631132
$ws.Range("N107").Value = -5290
# Row 134
$ws.Range("H134").Value = 33384368
$ws.Range("I134").Value = 47691110
$ws.Range("J134").Value = 1966.6666
$ws.Range("K134").Value = 143073330
$ws.Range("L134").Value = 5899.9998
$ws.Range("M134").Value = -143070795
$ws.Range("N134").Value = -10969.9998

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 15630934
$ws.Range("I31").Value = 18868926
$ws.Range("J31").Value = 29699.273
$ws.Range("K31").Value = 18868926
$ws.Range("L31").Value = 29699.273
$ws.Range("M31").Value = -18868631
$ws.Range("N31").Value = -30289.273
# Row 34
$ws.Range("H34").Value = 15630934
$ws.Range("I34").Value = 18868926
$ws.Range("J34").Value = 29699.273
$ws.Range("K34").Value = 18868926
$ws.Range("L34").Value = 29699.273
$ws.Range("M34").Value = -18868724
$ws.Range("N34").Value = -30103.273
# Row 50
$ws.Range("H50").Value = 10110.625
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 10110.625
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 10110.625
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -11360.625
# Row 51
$ws.Range("H51").Value = 10800.143
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 11120.2
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 11120.2
$ws.Range("M51").Value = -9264
$ws.Range("N51").Value = -12592.2
# Row 60
$ws.Range("H60").Value = 7355.6665
$ws.Range("I60").Value = 4997.5
$ws.Range("J60").Value = 8029.4287
$ws.Range("K60").Value = 4997.5
$ws.Range("L60").Value = 8029.4287
$ws.Range("M60").Value = -4486.5
$ws.Range("N60").Value = -9051.4287
# Row 61
$ws.Range("H61").Value = 10800.143
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 11120.2
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 11120.2
$ws.Range("M61").Value = -9652
$ws.Range("N61").Value = -11816.2
# Row 132
$ws.Range("H132").Value = 10758337
$ws.Range("I132").Value = 14495140
$ws.Range("J132").Value = 15026.5
$ws.Range("K132").Value = 43485420
$ws.Range("L132").Value = 45079.5
$ws.Range("M132").Value = -43482890
$ws.Range("N132").Value = -50139.5
# Row 134
$ws.Range("H134").Value = 26786874
$ws.Range("I134").Value = 33784908
$ws.Range("J134").Value = 5209608
$ws.Range("K134").Value = 101354724
$ws.Range("L134").Value = 15628824
$ws.Range("M134").Value = -101352189
$ws.Range("N134").Value = -15633894

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 16
$ws.Range("H16").Value = 751.9091
$ws.Range("I16").Value = 222.75
$ws.Range("J16").Value = 1054.2858
$ws.Range("K16").Value = 668.25
$ws.Range("L16").Value = 3162.8574
$ws.Range("M16").Value = -495.25
$ws.Range("N16").Value = -3508.8574

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1232.6666
$ws.Range("I113").Value = 1197.8
$ws.Range("J113").Value = 1319.8334
$ws.Range("K113").Value = 1197.8
$ws.Range("L113").Value = 1319.8334
$ws.Range("M113").Value = 972.2
$ws.Range("N113").Value = -5659.8334
# Row 122
$ws.Range("H122").Value = 125001200
$ws.Range("I122").Value = 166667440
$ws.Range("J122").Value = 2504
$ws.Range("K122").Value = 500002320
$ws.Range("L122").Value = 7512
$ws.Range("M122").Value = -499999870
$ws.Range("N122").Value = -12412
# Row 132
$ws.Range("H132").Value = 40003170
$ws.Range("I132").Value = 57144296
$ws.Range("J132").Value = 7203.6
$ws.Range("K132").Value = 171432888
$ws.Range("L132").Value = 21610.8
$ws.Range("M132").Value = -171430358
$ws.Range("N132").Value = -26670.8

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2172.182
$ws.Range("I7").Value = 2583.3333
$ws.Range("J7").Value = 1678.8
$ws.Range("K7").Value = 2583.3333
$ws.Range("L7").Value = 1678.8
$ws.Range("M7").Value = -2471.3333
$ws.Range("N7").Value = -1902.8
# Row 61
$ws.Range("H61").Value = 1392.6875
$ws.Range("I61").Value = 1021
$ws.Range("K61").Value = 1021
$ws.Range("M61").Value = -819
# Row 113
$ws.Range("H113").Value = 1392.6875
$ws.Range("I113").Value = 1021
$ws.Range("K113").Value = 1021
$ws.Range("M113").Value = 1149
# Row 126
$ws.Range("H126").Value = 2172.182
$ws.Range("I126").Value = 2583.3333
$ws.Range("J126").Value = 1678.8
$ws.Range("K126").Value = 7749.999899999999
$ws.Range("L126").Value = 5036.4
$ws.Range("M126").Value = -5279.999899999999
$ws.Range("N126").Value = -9976.4
# Row 132
$ws.Range("H132").Value = 4446047
$ws.Range("I132").Value = 6061757.5
$ws.Range("J132").Value = 2842.4167
$ws.Range("K132").Value = 18185272.5
$ws.Range("L132").Value = 8527.250100000001
$ws.Range("M132").Value = -18182742.5
$ws.Range("N132").Value = -13587.2501
# Row 136
$ws.Range("H136").Value = 3678.8655
$ws.Range("I136").Value = 3912.848
$ws.Range("J136").Value = 1885
$ws.Range("K136").Value = 11738.544
$ws.Range("L136").Value = 5655
$ws.Range("M136").Value = -9188.544
$ws.Range("N136").Value = -10755

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 5000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 5000
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -7746
# Row 107
$ws.Range("H107").Value = 45455010
$ws.Range("I107").Value = 542.5
$ws.Range("J107").Value = 166666910
$ws.Range("K107").Value = 1627.5
$ws.Range("L107").Value = 500000730
$ws.Range("M107").Value = 292.5
$ws.Range("N107").Value = -500004570
# Row 113
$ws.Range("H113").Value = 957.35
$ws.Range("I113").Value = 744.625
$ws.Range("J113").Value = 1099.1666
$ws.Range("K113").Value = 2233.875
$ws.Range("L113").Value = 3297.4998
$ws.Range("M113").Value = -63.875
$ws.Range("N113").Value = -7637.4998
# Row 126
$ws.Range("H126").Value = 208334320
$ws.Range("I126").Value = 83333970
$ws.Range("J126").Value = 333334660
$ws.Range("K126").Value = 250001910
$ws.Range("L126").Value = 1000003980
$ws.Range("M126").Value = -249999440
$ws.Range("N126").Value = -1000008920
# Row 132
$ws.Range("H132").Value = 7362100.5
$ws.Range("I132").Value = 4762923
$ws.Range("J132").Value = 9346927
$ws.Range("K132").Value = 14288769
$ws.Range("L132").Value = 28040781
$ws.Range("M132").Value = -14286239
$ws.Range("N132").Value = -28045841
# Row 136
$ws.Range("H136").Value = 12422313
$ws.Range("I136").Value = 6738142
$ws.Range("J136").Value = 31251132
$ws.Range("K136").Value = 20214426
$ws.Range("L136").Value = 93753396
$ws.Range("M136").Value = -20211876
$ws.Range("N136").Value = -93758496
